$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -6336
$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 8000
$ws.Range("K51").Value = 8000
$ws.Range("M51").Value = -7516
$ws.Range("H132").Value = 1638.1852
$ws.Range("I132").Value = 1547.3462
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4642.0386
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2112.0386
$ws.Range("N132").Value = -17060
$ws.Range("H137").Value = 2273.6667
$ws.Range("J137").Value = 3410.3333
$ws.Range("L137").Value = 10230.9999
$ws.Range("N137").Value = -15330.9999
$ws.Range("H141").Value = 1591.7693
$ws.Range("I141").Value = 1391.0834
$ws.Range("K141").Value = 4173.2502
$ws.Range("M141").Value = 1006.7498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14869.588
$ws.Range("I32").Value = 2594.5833
$ws.Range("J32").Value = 25780.703
$ws.Range("K32").Value = 2594.5833
$ws.Range("L32").Value = 25780.703
$ws.Range("M32").Value = -2307.5833
$ws.Range("N32").Value = -26354.703
$ws.Range("H61").Value = 3230.1667
$ws.Range("I61").Value = 3230.1667
$ws.Range("K61").Value = 3230.1667
$ws.Range("M61").Value = -3018.1667
$ws.Range("H132").Value = 2666
$ws.Range("I132").Value = 2230.6155
$ws.Range("J132").Value = 3798
$ws.Range("K132").Value = 6691.8465
$ws.Range("L132").Value = 11394
$ws.Range("M132").Value = -4161.8465
$ws.Range("N132").Value = -16454
$ws.Range("H136").Value = 3230.1667
$ws.Range("I136").Value = 3230.1667
$ws.Range("K136").Value = 9690.500100000001
$ws.Range("M136").Value = -7140.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4170.9165
$ws.Range("I105").Value = 3308.0527
$ws.Range("K105").Value = 3308.0527
$ws.Range("M105").Value = -1561.0527
$ws.Range("H107").Value = 771.6
$ws.Range("I107").Value = 778.64703
$ws.Range("K107").Value = 778.64703
$ws.Range("M107").Value = 1141.35297
$ws.Range("H134").Value = 2527.8845
$ws.Range("I134").Value = 2394.0952
$ws.Range("K134").Value = 7182.285600000001
$ws.Range("M134").Value = -4647.285600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3672
$ws.Range("I31").Value = 1840.4445
$ws.Range("J31").Value = 6968.8
$ws.Range("K31").Value = 1840.4445
$ws.Range("L31").Value = 6968.8
$ws.Range("M31").Value = -1545.4445
$ws.Range("N31").Value = -7558.8
$ws.Range("H34").Value = 3672
$ws.Range("I34").Value = 1840.4445
$ws.Range("J34").Value = 6968.8
$ws.Range("K34").Value = 1840.4445
$ws.Range("L34").Value = 6968.8
$ws.Range("M34").Value = -1638.4445
$ws.Range("N34").Value = -7372.8
$ws.Range("H58").Value = 3306.5
$ws.Range("I58").Value = 2612.6
$ws.Range("J58").Value = 3537.8
$ws.Range("K58").Value = 2612.6
$ws.Range("L58").Value = 3537.8
$ws.Range("M58").Value = -2409.6
$ws.Range("N58").Value = -3943.8
$ws.Range("H107").Value = 1105.909
$ws.Range("I107").Value = 910.6667
$ws.Range("K107").Value = 910.6667
$ws.Range("M107").Value = 1009.3333
$ws.Range("H132").Value = 3122.6667
$ws.Range("I132").Value = 2008.4117
$ws.Range("K132").Value = 6025.2351
$ws.Range("M132").Value = -3495.2351
$ws.Range("H134").Value = 3639.6875
$ws.Range("I134").Value = 2758.111
$ws.Range("J134").Value = 4773.143
$ws.Range("K134").Value = 8274.332999999999
$ws.Range("L134").Value = 14319.429
$ws.Range("M134").Value = -5739.332999999999
$ws.Range("N134").Value = -19389.429
$ws.Range("H136").Value = 3306.5
$ws.Range("I136").Value = 2612.6
$ws.Range("J136").Value = 3537.8
$ws.Range("K136").Value = 7837.799999999999
$ws.Range("L136").Value = 10613.4
$ws.Range("M136").Value = -5287.799999999999
$ws.Range("N136").Value = -15713.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H129").Value = 1738
$ws.Range("I129").Value = 1661
$ws.Range("K129").Value = 4983
$ws.Range("M129").Value = 17
$ws.Range("H132").Value = 1641.7778
$ws.Range("I132").Value = 1627
$ws.Range("J132").Value = 1693.5
$ws.Range("K132").Value = 14643
$ws.Range("L132").Value = 15241.5
$ws.Range("M132").Value = -12113
$ws.Range("N132").Value = -20301.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 45000
$ws.Range("J39").Value = 45000
$ws.Range("L39").Value = 45000
$ws.Range("N39").Value = -46064
$ws.Range("H132").Value = 2300.5264
$ws.Range("I132").Value = 1521.5333
$ws.Range("K132").Value = 4564.5999
$ws.Range("M132").Value = -2034.5999
$ws.Range("H141").Value = 58042.57
$ws.Range("J141").Value = 58042.57
$ws.Range("L141").Value = 58042.57
$ws.Range("N141").Value = -68402.57000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 14496.5
$ws.Range("J61").Value = 12499.5
$ws.Range("L61").Value = 12499.5
$ws.Range("N61").Value = -12903.5
$ws.Range("H113").Value = 14496.5
$ws.Range("J113").Value = 12499.5
$ws.Range("L113").Value = 12499.5
$ws.Range("N113").Value = -16839.5
$ws.Range("H132").Value = 4394.025
$ws.Range("I132").Value = 3806.4614
$ws.Range("J132").Value = 5485.2144
$ws.Range("K132").Value = 11419.3842
$ws.Range("L132").Value = 16455.6432
$ws.Range("M132").Value = -8889.3842
$ws.Range("N132").Value = -21515.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 879
$ws.Range("I81").Value = 879
$ws.Range("K81").Value = 1758
$ws.Range("M81").Value = -697
$ws.Range("H84").Value = 879
$ws.Range("I84").Value = 879
$ws.Range("K84").Value = 8790
$ws.Range("M84").Value = -3486
$ws.Range("H107").Value = 1182.8
$ws.Range("I107").Value = 532.5714
$ws.Range("K107").Value = 1597.7142
$ws.Range("M107").Value = 322.2857999999999
$ws.Range("H113").Value = 1117.6364
$ws.Range("J113").Value = 1122.8572
$ws.Range("L113").Value = 3368.5716
$ws.Range("N113").Value = -7708.571599999999
$ws.Range("H122").Value = 839.25
$ws.Range("I122").Value = 839.25
$ws.Range("K122").Value = 2517.75
$ws.Range("M122").Value = -67.75
$ws.Range("H136").Value = 3170.6316
$ws.Range("I136").Value = 1391
$ws.Range("K136").Value = 4173
$ws.Range("M136").Value = -1623

Write-Output "Applied all cell updates."